# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Column G ("K") values are recalculated (std/mean regen + s_vals calc) and
# rewritten for each data row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$kValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 3
    12 = 3
    13 = 2
    14 = 1
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 2
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 0
    26 = 3
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 1
    32 = 3
    33 = 1
    34 = 1
    35 = 0
    36 = 2
    37 = 1
    38 = 2
    40 = 1
    41 = 2
    42 = 0
    43 = 3
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    48 = 1
    49 = 2
    50 = 0
    51 = 0
    52 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
